$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 16 with the latest quarter data (01-07-2021).
# Force the date-like label to be stored as text (matching column A's
# existing text entries) instead of being auto-converted to a date serial.
$ws.Cells.Item(16, 1).NumberFormat = "@"
$ws.Cells.Item(16, 1).Value = "01-07-2021"
$ws.Cells.Item(16, 1).Style = "Normal"

$ws.Cells.Item(16, 2).Value = 115.56
$ws.Cells.Item(16, 3).Value = 108.79
$ws.Cells.Item(16, 4).Value = 99.16
$ws.Cells.Item(16, 5).Value = 112
$ws.Cells.Item(16, 6).Value = 111.89
$ws.Cells.Item(16, 7).Value = 107.28
$ws.Cells.Item(16, 8).Value = 112.12
$ws.Cells.Item(16, 9).Value = 93.44
$ws.Cells.Item(16, 10).Value = 109.28
$ws.Cells.Item(16, 11).Value = 112.24
$ws.Cells.Item(16, 12).Value = 110.65
$ws.Cells.Item(16, 13).Value = 111.75
